$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '41.771.01'
$ws.Range('E2').Value = '  -3.71%  '
$ws.Range('D3').Value = '2.238.91'
$ws.Range('E3').Value = '  -4.05%  '
$ws.Range('E4').Value = '  -0.17%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '249.05'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +4.17%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.632'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -4.15%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '71.71'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -2.28%  '
$ws.Range('E8').Value = '  +0.05%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.561'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -4.81%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '40.84'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +13.12%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0968'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -5.54%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '58.67'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -1.06%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.105'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -2.41%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.88'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -5.02%  '
$ws.Range('D15').Value = '2.573.72'
$ws.Range('E15').Value = '  -3.99%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '15.02'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -7.23%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.854'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -7.56%  '
$ws.Range('D18').Value = '2.235.98'
$ws.Range('E18').Value = '  -4.15%  '
$ws.Range('D19').Value = '41.732.22'
$ws.Range('E19').Value = '  -3.71%  '
$ws.Range('D20').Value = '0.0₃0960'
$ws.Range('E20').Value = '  -6.95%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '73.33'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -4.36%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.16'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -5.44%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '235.24'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -6.22%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.23'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +20.23%  '
$ws.Range('E25').Value = '  -0.02%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '3.72'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.78%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.47'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.31%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.05'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -4.71%  '
$ws.Range('E29').Value = '  -3.93%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '171.06'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -1.48%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '20.73'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -5.37%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.122'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -4.22%  '
$ws.Range('E33').Value = '  -6.06%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0722'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -3.02%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.34'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -3.99%  '
$ws.Range('B36').Value = 'InjectiveProtocol'
$ws.Range('C36').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '26.60'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +26.58%  '
$ws.Range('B37').Value = 'Filecoin'
$ws.Range('C37').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '4.70'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -7.98%  '
$ws.Range('B38').Value = 'RenderToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '4.11'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +9.46%  '
$ws.Range('E39').Value = '  +0.22%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.28'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -3.47%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '6.01'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -7.33%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '66.69'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.90%  '
$ws.Range('B43').Value = 'Algorand'
$ws.Range('C43').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.212'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +6.47%  '
$ws.Range('B44').Value = 'FTXToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '5.12'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -2.72%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '11.74'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +19.76%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '8.80'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -3.40%  '
$ws.Range('B47').Value = 'SynthetixNetwork'
$ws.Range('C47').Value = 'https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '4.78'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +6.71%  '
$ws.Range('B48').Value = 'Cronos'
$ws.Range('C48').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.101'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -5.94%  '
$ws.Range('E49').Value = '  +0.03%  '
$ws.Range('B50').Value = 'BitTorrent-New'
$ws.Range('C50').Value = 'https://coinranking.com/coin/w4MqH_Xe8+bittorrent-new-btt'
$ws.Range('D50').Value = '0.0₃0153'
$ws.Range('E50').Value = '  +13.97%  '
$ws.Range('B51').Value = 'TrustWalletToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.18'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -3.78%  '
